$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $addr, $val)
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

Set-TextValue $ws 'D2' '29.436.01'
Set-TextValue $ws 'E2' '  +0.13%  '
Set-TextValue $ws 'D3' '1.848.65'
Set-TextValue $ws 'E3' '  +0.26%  '
Set-TextValue $ws 'E4' '  +0.22%  '
Set-TextValue $ws 'D5' '240.78'
Set-TextValue $ws 'E5' '  +0.74%  '
Set-TextValue $ws 'D6' '0.6262'
Set-TextValue $ws 'E6' '  -0.80%  '
Set-TextValue $ws 'E7' '  +0.17%  '
Set-TextValue $ws 'D8' '0.07673'
Set-TextValue $ws 'E8' '  +1.86%  '
Set-TextValue $ws 'D9' '0.2913'
Set-TextValue $ws 'E9' '  -0.59%  '
Set-TextValue $ws 'D10' '24.74'
Set-TextValue $ws 'E10' '  +1.21%  '
Set-TextValue $ws 'D11' '0.07749'
Set-TextValue $ws 'E11' '  +0.49%  '
Set-TextValue $ws 'D12' '1.843.41'
Set-TextValue $ws 'E12' '  -0.51%  '
Set-TextValue $ws 'D13' '5.026'
Set-TextValue $ws 'E13' '  +0.57%  '
Set-TextValue $ws 'D14' '0.6811'
Set-TextValue $ws 'E14' '  +0.32%  '
Set-TextValue $ws 'D15' '0.00001072'
Set-TextValue $ws 'D16' '83.50'
Set-TextValue $ws 'E16' '  +0.57%  '
Set-TextValue $ws 'D17' '6.165'
Set-TextValue $ws 'E17' '  +0.03%  '
Set-TextValue $ws 'D18' '29.462.76'
Set-TextValue $ws 'E18' '  +0.12%  '
Set-TextValue $ws 'D19' '228.12'
Set-TextValue $ws 'E19' '  +0.05%  '
Set-TextValue $ws 'E20' '  -0.23%  '
Set-TextValue $ws 'D21' '1.002'
Set-TextValue $ws 'E21' '  +0.18%  '
Set-TextValue $ws 'D22' '7.414'
Set-TextValue $ws 'E22' '  -0.57%  '
Set-TextValue $ws 'E23' '  +0.14%  '
Set-TextValue $ws 'D24' '157.48'
Set-TextValue $ws 'E24' '  +0.30%  '
Set-TextValue $ws 'E25' '  -1.39%  '
Set-TextValue $ws 'D26' '8.394'
Set-TextValue $ws 'E26' '  +0.30%  '
Set-TextValue $ws 'D27' '17.68'
Set-TextValue $ws 'E27' '  +0.60%  '
Set-TextValue $ws 'D28' '1.343'
Set-TextValue $ws 'E28' '  +5.13%  '
Set-TextValue $ws 'D29' '1.466'
Set-TextValue $ws 'E29' '  +0.49%  '
Set-TextValue $ws 'D30' '0.05629'
Set-TextValue $ws 'E31' '  +0.25%  '
Set-TextValue $ws 'E32' '  +0.04%  '
Set-TextValue $ws 'E33' '  +0.37%  '
Set-TextValue $ws 'E34' '  +0.33%  '
Set-TextValue $ws 'D35' '0.7079'
Set-TextValue $ws 'E35' '  -0.52%  '
Set-TextValue $ws 'E36' '  +0.33%  '
Set-TextValue $ws 'D37' '1.230.39'
Set-TextValue $ws 'E37' '  -1.24%  '
Set-TextValue $ws 'D38' '2.769'
Set-TextValue $ws 'E38' '  +0.24%  '
Set-TextValue $ws 'D39' '0.01786'
Set-TextValue $ws 'E39' '  -1.22%  '
Set-TextValue $ws 'D40' '6.539'
Set-TextValue $ws 'E40' '  +3.49%  '
Set-TextValue $ws 'D41' '0.9041'
Set-TextValue $ws 'B43' 'RocketPoolETH'
Set-TextValue $ws 'C43' 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
Set-TextValue $ws 'D43' '2.033.51'
Set-TextValue $ws 'E43' '  -0.54%  '
Set-TextValue $ws 'B44' 'Quant'
Set-TextValue $ws 'C44' 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue $ws 'D44' '101.75'
Set-TextValue $ws 'E44' '  -0.04%  '
Set-TextValue $ws 'B45' 'Aave'
Set-TextValue $ws 'C45' 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws 'D45' '65.84'
Set-TextValue $ws 'E45' '  +0.10%  '
Set-TextValue $ws 'B46' 'BabyDogeCoin'
Set-TextValue $ws 'C46' 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue $ws 'D46' '0.00000000120'
Set-TextValue $ws 'B47' 'Aptos'
Set-TextValue $ws 'C47' 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws 'D47' '7.153'
Set-TextValue $ws 'E47' '  +0.91%  '
Set-TextValue $ws 'B48' 'TheSandbox'
Set-TextValue $ws 'C48' 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue $ws 'D48' '0.4012'
Set-TextValue $ws 'E48' '  +0.28%  '
Set-TextValue $ws 'D49' '0.1151'
Set-TextValue $ws 'E49' '  +2.78%  '
Set-TextValue $ws 'B50' 'EnergySwap'
Set-TextValue $ws 'C50' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws 'D50' '8.973'
Set-TextValue $ws 'E50' '  +0.96%  '
Set-TextValue $ws 'B51' 'RenderToken'
Set-TextValue $ws 'C51' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws 'D51' '1.672'
Set-TextValue $ws 'E51' '  +0.08%  '
